$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Find the last used row based on column A (Beteckning)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C = "Förändrad" date column. Every data row (2..lastRow) had its
# value bumped by exactly one day (45178 -> 45179).
$rng = $ws.Range("C2:C$lastRow")
$rng.Value = 45179
